$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '64.471.43'
$ws.Range("E2").Value = '  +0.10%  '
$ws.Range("D3").Value = '3.145.86'
$ws.Range("E3").Value = '  -0.32%  '
$ws.Range("E4").Value = '  +0.10%  '
$ws.Range("D5").Value = '611.14'
$ws.Range("E5").Value = '  +0.46%  '
$ws.Range("D6").Value = '143.78'
$ws.Range("E6").Value = '  -2.04%  '
$ws.Range("E7").Value = '  +0.04%  '
$ws.Range("D8").Value = '3.146.18'
$ws.Range("E8").Value = '  -0.18%  '
$ws.Range("D9").Value = '0.527'
$ws.Range("E9").Value = '  -0.02%  '
$ws.Range("E10").Value = '  -0.14%  '
$ws.Range("E11").Value = '  -3.18%  '
$ws.Range("D12").Value = '0.473'
$ws.Range("E12").Value = '  -0.42%  '
$ws.Range("E13").Value = '  +0.83%  '
$ws.Range("D14").Value = '35.52'
$ws.Range("E14").Value = '  -1.51%  '
$ws.Range("D15").Value = '3.663.92'
$ws.Range("E15").Value = '  -0.19%  '
$ws.Range("E16").Value = '  +2.95%  '
$ws.Range("D17").Value = '64.409.12'
$ws.Range("E17").Value = '  +0.07%  '
$ws.Range("D18").Value = '3.172.33'
$ws.Range("E18").Value = '  +0.53%  '
$ws.Range("D19").Value = '6.85'
$ws.Range("E19").Value = '  -1.50%  '
$ws.Range("D20").Value = '476.29'
$ws.Range("E20").Value = '  -0.42%  '
$ws.Range("D21").Value = '14.67'
$ws.Range("E21").Value = '  +0.50%  '
$ws.Range("D22").Value = '0.722'
$ws.Range("E22").Value = '  +1.68%  '
$ws.Range("D23").Value = '7.83'
$ws.Range("E23").Value = '  +1.01%  '
$ws.Range("D24").Value = '13.67'
$ws.Range("E24").Value = '  -0.76%  '
$ws.Range("D25").Value = '84.86'
$ws.Range("E25").Value = '  +1.61%  '
$ws.Range("D26").Value = '0.999'
$ws.Range("E26").Value = '  -0.09%  '
$ws.Range("D27").Value = '2.79'
$ws.Range("E27").Value = '  -3.30%  '
$ws.Range("D28").Value = '8.59'
$ws.Range("E28").Value = '  +1.99%  '
$ws.Range("D29").Value = '7.40'
$ws.Range("E29").Value = '  +8.27%  '
$ws.Range("E30").Value = '  +0.79%  '
$ws.Range("D31").Value = '2.09'
$ws.Range("E31").Value = '  -4.68%  '
$ws.Range("E32").Value = '  +0.02%  '
$ws.Range("D33").Value = '26.62'
$ws.Range("E33").Value = '  +1.59%  '
$ws.Range("D34").Value = '2.65'
$ws.Range("E34").Value = '  -3.99%  '
$ws.Range("E35").Value = '  +1.13%  '
$ws.Range("D36").Value = '5.95'
$ws.Range("E36").Value = '  -1.24%  '
$ws.Range("D37").Value = '52.59'
$ws.Range("E37").Value = '  -2.88%  '
$ws.Range("D38").Value = '0.0₃0739'
$ws.Range("E38").Value = '  +2.70%  '
$ws.Range("D39").Value = '3.02'
$ws.Range("E39").Value = '  +3.45%  '
$ws.Range("D40").Value = '453.22'
$ws.Range("E40").Value = '  +0.13%  '
$ws.Range("D41").Value = '0.0395'
$ws.Range("E41").Value = '  -0.51%  '
$ws.Range("E42").Value = '  -0.37%  '
$ws.Range("D43").Value = '8.33'
$ws.Range("E43").Value = '  -1.44%  '
$ws.Range("D44").Value = '2.855.64'
$ws.Range("E44").Value = '  +0.40%  '
$ws.Range("D45").Value = '0.267'
$ws.Range("E45").Value = '  -0.63%  '
$ws.Range("E46").Value = '  -0.02%  '
$ws.Range("D47").Value = '2.43'
$ws.Range("E47").Value = '  +4.82%  '
$ws.Range("D48").Value = '26.44'
$ws.Range("E48").Value = '  +0.01%  '
$ws.Range("E50").Value = '  -0.19%  '
$ws.Range("D51").Value = '120.09'
$ws.Range("E51").Value = '  +1.30%  '
